$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet1"

$src = $ws1.Range("A9:K37")
$src.Copy()
$dst = $newSheet.Range("A1")
$dst.PasteSpecial()
